$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Text changes (shared strings content) ---
# "Remarks" -> "Test result" (header cell E4 on both sheets)
$ws1.Range("E4").Value = "Test result"
$ws2.Range("E4").Value = "Test result"

# "Ok" -> "Success" (result cells)
$ws1.Range("E5:E8").Value = "Success"
$ws2.Range("E5:E10").Value = "Success"

# --- Column widths (new column E sizing) ---
$ws1.Columns.Item(5).ColumnWidth = 10.333333
$ws2.Columns.Item(5).ColumnWidth = 19.333333

# --- Sheet view / selection changes ---
# Sheet2: selection moves from A4:E10 to E4:E10, and it stops being the active/tab-selected sheet.
[void]$ws2.Select()
[void]$ws2.Range("E4:E10").Select()

# Sheet1: becomes the active/tab-selected sheet, with D19 selected.
# (Selected last so it ends up as the workbook's active tab.)
[void]$ws1.Select()
[void]$ws1.Range("D19").Select()
